$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Drop the rows that held MedicationAgreement / AdministrationAgreement
#        mappings in column B (old rows 26-56 of the tail, the others in
#        20-25 get overwritten in-place below with the MedicationUse rows
#        that used to start further down). ---
$ws.Range("A26:B56").EntireRow.Delete()

# --- 2. Rewrite column B with the MedicationUse.* mapping values, reusing
#        the existing body style (same as column A) for every touched cell.
#        Copy/PasteSpecial(formats) keeps reusing style index 2 rather than
#        minting a duplicate cellXf. ---
$mapping = @{
    2  = "MedicationUse"
    13 = "MedicationUse.UseIndicator"
    14 = "MedicationUse.ProductUsed::PharmaceuticalProduct"
    16 = "MedicationUse.ReasonForUse"
    18 = "MedicationUse.InstructionsForUse"
    19 = "MedicationUse.PeriodOfUse::TimeInterval"
    20 = "MedicationUse.Prescriber::HealthProfessional"
    21 = "MedicationUse.MedicationUseDateTime"
    22 = "MedicationUse.AsAgreedIndicator"
    23 = "MedicationUse.MedicationUseStopType"
    24 = "MedicationUse.ReasonModificationOrDiscontinuationOfUse"
    25 = "MedicationUse.Comment"
}

$ws.Range("A3").Copy()
foreach ($row in $mapping.Keys) {
    $cell = $ws.Range("B$row")
    $cell.Value = $mapping[$row]
    $cell.PasteSpecial(-4122)
}
$excel.CutCopyMode = $false

# --- 3. Row 21 used to carry a custom 15pt height (leftover from the old
#        content); the new content no longer needs it. ---
$ws.Rows.Item(21).AutoFit()

# --- 4. Match the saved selection / view state. ---
$ws.Range("A23:A24").EntireRow.Select()
